$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out existing data (A1:B3) since the new layout only uses column A (A1:A6)
$ws.Range("A1:B3").Clear()

# Set the new values: A1:A6 all contain "DineshQA"
$ws.Range("A1").Value = "DineshQA"
$ws.Range("A2").Value = "DineshQA"
$ws.Range("A3").Value = "DineshQA"
$ws.Range("A4").Value = "DineshQA"
$ws.Range("A5").Value = "DineshQA"
$ws.Range("A6").Value = "DineshQA"
